$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data rows (2-7) and extend with new data (rows 2-10)
$ws.Range("A2:T7").ClearContents()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Agt"
$ws.Range("C2").Value = "Agtr1a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1206283333333333
$ws.Range("H2").Value = 0.361885
$ws.Range("I2").Value = 0.09993369694616584
$ws.Range("J2").Value = 0.09993369694616584
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1356863333333333
$ws.Range("N2").Value = 0.407059
$ws.Range("O2").Value = 0.004454204096299941
$ws.Range("P2").Value = 0.004454204096299941
$ws.Range("Q2").Value = 0.01636761624611111
$ws.Range("R2").Value = 0.147308546215
$ws.Range("S2").Value = 0.0004451250822960088
$ws.Range("T2").Value = 0.0004451250822960088

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Agt"
$ws.Range("C3").Value = "Agtr1a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1206283333333333
$ws.Range("H3").Value = 0.361885
$ws.Range("I3").Value = 0.09993369694616584
$ws.Range("J3").Value = 0.09993369694616584
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 17.195945
$ws.Range("N3").Value = 51.587835
$ws.Range("O3").Value = 0.564494940478519
$ws.Range("P3").Value = 0.5644949404785189
$ws.Range("Q3").Value = 2.074318185441666
$ws.Range("R3").Value = 18.668863668975
$ws.Range("S3").Value = 0.05641206630942425
$ws.Range("T3").Value = 0.05641206630942423

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Agt"
$ws.Range("C4").Value = "Agtr1a"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1206283333333333
$ws.Range("H4").Value = 0.361885
$ws.Range("I4").Value = 0.09993369694616584
$ws.Range("J4").Value = 0.09993369694616584
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 13.13090033333333
$ws.Range("N4").Value = 39.392701
$ws.Range("O4").Value = 0.4310508554251812
$ws.Range("P4").Value = 0.4310508554251811
$ws.Range("Q4").Value = 1.583958622376111
$ws.Range("R4").Value = 14.255627601385
$ws.Range("S4").Value = 0.0430765055544456
$ws.Range("T4").Value = 0.04307650555444559

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Agt"
$ws.Range("C5").Value = "Agtr1a"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4683593333333333
$ws.Range("H5").Value = 1.405078
$ws.Range("I5").Value = 0.3880090057275787
$ws.Range("J5").Value = 0.3880090057275787
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1356863333333333
$ws.Range("N5").Value = 0.407059
$ws.Range("O5").Value = 0.004454204096299941
$ws.Range("P5").Value = 0.004454204096299941
$ws.Range("Q5").Value = 0.06354996062244445
$ws.Range("R5").Value = 0.571949645602
$ws.Range("S5").Value = 0.001728271302713048
$ws.Range("T5").Value = 0.001728271302713048

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Agt"
$ws.Range("C6").Value = "Agtr1a"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.4683593333333333
$ws.Range("H6").Value = 1.405078
$ws.Range("I6").Value = 0.3880090057275787
$ws.Range("J6").Value = 0.3880090057275787
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 17.195945
$ws.Range("N6").Value = 51.587835
$ws.Range("O6").Value = 0.564494940478519
$ws.Range("P6").Value = 0.5644949404785189
$ws.Range("Q6").Value = 8.053881336236666
$ws.Range("R6").Value = 72.48493202613
$ws.Range("S6").Value = 0.2190291205933189
$ws.Range("T6").Value = 0.2190291205933188

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Agt"
$ws.Range("C7").Value = "Agtr1a"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.4683593333333333
$ws.Range("H7").Value = 1.405078
$ws.Range("I7").Value = 0.3880090057275787
$ws.Range("J7").Value = 0.3880090057275787
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 13.13090033333333
$ws.Range("N7").Value = 39.392701
$ws.Range("O7").Value = 0.4310508554251812
$ws.Range("P7").Value = 0.4310508554251811
$ws.Range("Q7").Value = 6.149979726186445
$ws.Range("R7").Value = 55.34981753567801
$ws.Range("S7").Value = 0.1672516138315468
$ws.Range("T7").Value = 0.1672516138315468

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Agt"
$ws.Range("C8").Value = "Agtr1a"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.618096
$ws.Range("H8").Value = 1.854288
$ws.Range("I8").Value = 0.5120572973262555
$ws.Range("J8").Value = 0.5120572973262555
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1356863333333333
$ws.Range("N8").Value = 0.407059
$ws.Range("O8").Value = 0.004454204096299941
$ws.Range("P8").Value = 0.004454204096299941
$ws.Range("Q8").Value = 0.08386717988799999
$ws.Range("R8").Value = 0.754804618992
$ws.Range("S8").Value = 0.002280807711290884
$ws.Range("T8").Value = 0.002280807711290884

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Agt"
$ws.Range("C9").Value = "Agtr1a"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.618096
$ws.Range("H9").Value = 1.854288
$ws.Range("I9").Value = 0.5120572973262555
$ws.Range("J9").Value = 0.5120572973262555
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 17.195945
$ws.Range("N9").Value = 51.587835
$ws.Range("O9").Value = 0.564494940478519
$ws.Range("P9").Value = 0.5644949404785189
$ws.Range("Q9").Value = 10.62874482072
$ws.Range("R9").Value = 95.65870338648
$ws.Range("S9").Value = 0.2890537535757759
$ws.Range("T9").Value = 0.2890537535757758

# Row 10
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Agt"
$ws.Range("C10").Value = "Agtr1a"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.618096
$ws.Range("H10").Value = 1.854288
$ws.Range("I10").Value = 0.5120572973262555
$ws.Range("J10").Value = 0.5120572973262555
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 13.13090033333333
$ws.Range("N10").Value = 39.392701
$ws.Range("O10").Value = 0.4310508554251812
$ws.Range("P10").Value = 0.4310508554251811
$ws.Range("Q10").Value = 8.116156972432
$ws.Range("R10").Value = 73.04541275188801
$ws.Range("S10").Value = 0.2207227360391888
$ws.Range("T10").Value = 0.2207227360391887
